# Applies the cryptos.xlsx price/volume refresh described in the commit diff.
# (GitHub Actions scheduled data refresh: row edits + two name/link/price/volume swaps.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cells whose new text would otherwise be auto-parsed as a number by Excel
# (e.g. "9.70" -> 9.7). Force the Text number format just long enough to assign
# the literal string, then restore the default style so no formatting is left
# behind (matches the source cells, which carry no style attribute).
$textCoercedCells = @(
    'D5',
    'D8',
    'D10',
    'D12',
    'D14',
    'D15',
    'D16',
    'D19',
    'D20',
    'D22',
    'D24',
    'D25',
    'D26',
    'D27',
    'D28',
    'D29',
    'D32',
    'D33',
    'D34',
    'D37',
    'D39',
    'D40',
    'D41',
    'D42',
    'D45',
    'D49',
    'D50',
)
foreach ($addr in $textCoercedCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# --- Plain value updates (row 2-22, 26-51 price/volume refresh)
$ws.Range('D2').Value = '37.593.47'
$ws.Range('E2').Value = '  -0.38%  '
$ws.Range('D3').Value = '2.079.95'
$ws.Range('E3').Value = '  +0.14%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('E5').Value = '  -0.09%  '
$ws.Range('E6').Value = '  +2.16%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('E8').Value = '  -0.32%  '
$ws.Range('E9').Value = '  -0.33%  '
$ws.Range('E10').Value = '  -0.60%  '
$ws.Range('E11').Value = '  +3.11%  '
$ws.Range('E12').Value = '  +2.11%  '
$ws.Range('D13').Value = '2.384.72'
$ws.Range('E13').Value = '  +0.06%  '
$ws.Range('E14').Value = '  +1.37%  '
$ws.Range('E15').Value = '  -0.02%  '
$ws.Range('E16').Value = '  +1.07%  '
$ws.Range('D17').Value = '2.081.88'
$ws.Range('E17').Value = '  +0.26%  '
$ws.Range('D18').Value = '37.560.94'
$ws.Range('E18').Value = '  -0.34%  '
$ws.Range('E19').Value = '  -2.03%  '
$ws.Range('E20').Value = '  -0.48%  '
$ws.Range('D21').Value = '0.0₃0835'
$ws.Range('E21').Value = '  +0.09%  '
$ws.Range('E22').Value = '  +0.16%  '
$ws.Range('B24').Value = 'PancakeSwap'
$ws.Range('C24').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('E24').Value = '  -0.73%  '
$ws.Range('B25').Value = 'Toncoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('E25').Value = '  -2.62%  '
$ws.Range('E26').Value = '  +7.38%  '
$ws.Range('E27').Value = '  +0.38%  '
$ws.Range('E28').Value = '  -3.96%  '
$ws.Range('E29').Value = '  +0.16%  '
$ws.Range('E30').Value = '  -0.09%  '
$ws.Range('E31').Value = '  +0.98%  '
$ws.Range('E32').Value = '  -0.08%  '
$ws.Range('E33').Value = '  +1.18%  '
$ws.Range('E34').Value = '  -0.21%  '
$ws.Range('E35').Value = '  +1.07%  '
$ws.Range('E36').Value = '  -0.68%  '
$ws.Range('E37').Value = '  -1.84%  '
$ws.Range('E38').Value = '  -0.09%  '
$ws.Range('E39').Value = '  +2.15%  '
$ws.Range('E40').Value = '  +8.58%  '
$ws.Range('E41').Value = '  +2.86%  '
$ws.Range('E42').Value = '  -1.51%  '
$ws.Range('E43').Value = '  +3.94%  '
$ws.Range('E44').Value = '  +0.55%  '
$ws.Range('B45').Value = 'InjectiveProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('E45').Value = '  +3.64%  '
$ws.Range('B46').Value = 'Maker'
$ws.Range('C46').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D46').Value = '1.465.85'
$ws.Range('E46').Value = '  +1.04%  '
$ws.Range('E47').Value = '  -1.12%  '
$ws.Range('E48').Value = '  -5.54%  '
$ws.Range('E49').Value = '  -2.25%  '
$ws.Range('E50').Value = '  -2.23%  '
$ws.Range('D51').Value = '2.268.90'
$ws.Range('E51').Value = '  +0.02%  '

# --- Text-coerced value updates
$ws.Range('D5').Value = '233.53'
$ws.Range('D8').Value = '58.09'
$ws.Range('D10').Value = '0.0779'
$ws.Range('D12').Value = '15.16'
$ws.Range('D14').Value = '21.15'
$ws.Range('D15').Value = '0.773'
$ws.Range('D16').Value = '5.36'
$ws.Range('D19').Value = '6.05'
$ws.Range('D20').Value = '70.78'
$ws.Range('D22').Value = '229.18'
$ws.Range('D24').Value = '2.38'
$ws.Range('D25').Value = '2.34'
$ws.Range('D26').Value = '9.70'
$ws.Range('D27').Value = '170.60'
$ws.Range('D28').Value = '0.133'
$ws.Range('D29').Value = '19.48'
$ws.Range('D32').Value = '4.66'
$ws.Range('D33').Value = '0.0638'
$ws.Range('D34').Value = '4.64'
$ws.Range('D37').Value = '3.33'
$ws.Range('D39').Value = '5.36'
$ws.Range('D40').Value = '0.0233'
$ws.Range('D41').Value = '100.82'
$ws.Range('D42').Value = '0.0957'
$ws.Range('D45').Value = '16.93'
$ws.Range('D49').Value = '7.22'
$ws.Range('D50').Value = '2.94'

# Restore default (no explicit) style on the coerced cells.
foreach ($addr in $textCoercedCells) {
    $ws.Range($addr).Style = "Normal"
}
